$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data taken from the target diff (rows 4-10, columns A-I).
# "I" holds a dd/mm/yyyy-looking string that must stay literal TEXT,
# not get auto-converted into an Excel date serial number.
$rowsData = @(
    @{ A=2; B='5010754-58.2017.8.21.0001'; C='0196807-38.2017.8.21.0001'; D='Digitalizado'; E='Sem dados de processo originário 2'; F='Nulo'; G='Sem dados de processo originário 3'; H='Nulo'; I='10/11/2017' },
    @{ A=3; B='5008281-41.2013.8.21.0001'; C='0087776-25.2013.8.21.0001'; D='Digitalizado'; E='Sem dados de processo originário 2'; F='Nulo'; G='Sem dados de processo originário 3'; H='Nulo'; I='26/12/2011' },
    @{ A=4; B='5001221-76.2018.8.21.0054'; C='0066578-85.2019.8.21.9000'; D='Relacionado na TR'; E='9000149-83.2018.8.21.0054'; F='Migrado'; G='Sem dados de processo originário 3'; H='Nulo'; I='13/03/2018' },
    @{ A=5; B='5006432-89.2018.8.21.0023'; C='9003986-45.2018.8.21.0023'; D='Migrado'; E='Sem dados de processo originário 2'; F='Nulo'; G='Sem dados de processo originário 3'; H='Nulo'; I='31/10/2018' },
    @{ A=6; B='5006429-37.2018.8.21.0023'; C='9003787-23.2018.8.21.0023'; D='Migrado'; E='Sem dados de processo originário 2'; F='Nulo'; G='Sem dados de processo originário 3'; H='Nulo'; I='17/10/2018' },
    @{ A=7; B='5009614-96.2011.8.21.0001'; C='0420415-91.2011.8.21.0001'; D='Digitalizado'; E='Sem dados de processo originário 2'; F='Nulo'; G='Sem dados de processo originário 3'; H='Nulo'; I='07/12/2011' },
    @{ A=8; B='5033806-25.2013.8.21.0001'; C='0007850-16.2013.8.21.3001'; D='Digitalizado'; E='Sem dados de processo originário 2'; F='Nulo'; G='Sem dados de processo originário 3'; H='Nulo'; I='04/03/2013' }
)

# Scratch cell well outside the used range. Routing the date-like
# strings through a formula -> Copy -> PasteSpecial(xlPasteValues)
# round trip lands them as literal text instead of triggering the
# automatic "looks like a date" parsing that a direct .Value assignment
# of a bare dd/mm/yyyy string would trigger.
$scratch = $ws.Cells.Item(1, 50)

$r = 4
foreach ($row in $rowsData) {
    # A: numeric index, styled like the existing A4 cell (border/center/bold/top)
    $ws.Cells.Item(4, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Cells.Item($r, 1).Value = $row.A

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H

    # I: force literal text so the dd/mm/yyyy string is not parsed as a date
    $scratch.Formula = '="' + $row.I + '"'
    $scratch.Copy() | Out-Null
    $ws.Cells.Item($r, 9).PasteSpecial(-4163) | Out-Null   # xlPasteValues

    $r = $r + 1
}

$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false
